$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 731, shifting the existing rows 731-760 down to 732-761.
$ws.Rows("731:731").Insert()

# Populate the newly inserted row 731 with the new record's data.
$ws.Cells.Item(731, 1).Value2 = 5
$ws.Cells.Item(731, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(731, 3).Value2 = "Maule"
$ws.Cells.Item(731, 4).Value2 = 44939
$ws.Cells.Item(731, 5).Value2 = 7
$ws.Cells.Item(731, 6).Value2 = "Fruta"
$ws.Cells.Item(731, 7).Value2 = 100101
$ws.Cells.Item(731, 8).Value2 = "Berries"
$ws.Cells.Item(731, 9).Value2 = 100112025
$ws.Cells.Item(731, 10).Value2 = "Frutilla"
$ws.Cells.Item(731, 11).Value2 = "Sin especificar"
$ws.Cells.Item(731, 12).Value2 = "Primera"
$ws.Cells.Item(731, 13).Value2 = 250
$ws.Cells.Item(731, 14).Value2 = 8000
$ws.Cells.Item(731, 15).Value2 = 8000
$ws.Cells.Item(731, 16).Value2 = 8000
$ws.Cells.Item(731, 17).Value2 = "`$/caja 7 kilos"
$ws.Cells.Item(731, 18).Value2 = "Región del Maule"
$ws.Cells.Item(731, 19).Value2 = 1143
$ws.Cells.Item(731, 20).Value2 = 7
